$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the key (A-column) entries for the building-block pairs that were
# retired in this revision of the "xyz-part of 32111" key/value list.
# Clearing (rather than deleting rows) lets Excel's own shared-string
# compaction + empty-cell/row pruning do the rest, matching how this was
# produced originally.
$rangesToClear = @(
    "A80:A84",
    "A112:A117",
    "A230",
    "A238",
    "A251",
    "A262",
    "A270",
    "A283",
    "A304",
    "A320"
)

foreach ($addr in $rangesToClear) {
    $ws.Range($addr).ClearContents()
}

# Leave the UI selection where the author's session ended up.
$ws.Range("A321").Select()
